# Edit script: insert 6 new price rows for "Vega Monumental Concepción - Melón"
# right after the existing block for date 2022-06-18 worth of rows (old row 266),
# which shifts all subsequent rows down by 6 (old row 267 -> new row 273, etc.),
# and correspondingly appends the tail (old rows 339-344) as new rows 345-350.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 6 blank rows before the old row 267 (shifts everything down by 6 rows)
$ws.Range("A267:A272").EntireRow.Insert()

# 2) Fill the 6 new rows with the new price data (market report date 2023-02-07)
$newDate = Get-Date -Year 2023 -Month 2 -Day 7 -Hour 0 -Minute 0 -Second 0

$newRows = @(
    @("Calameño", "Extra",   1000, 1100, 1100, 1100, "Región de O'Higgins", 1100),
    @("Calameño", "Primera", 1000, 900,  900,  900,  "Región de O'Higgins", 900),
    @("Calameño", "Segunda", 500,  700,  700,  700,  "Región de O'Higgins", 700),
    @("Tuna",     "Extra",   1000, 1100, 1100, 1100, "Región de O'Higgins", 1100),
    @("Tuna",     "Primera", 1000, 900,  900,  900,  "Región de O'Higgins", 900),
    @("Tuna",     "Segunda", 500,  700,  700,  700,  "Región de O'Higgins", 700)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = 267 + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = 11
    $ws.Cells.Item($r, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($r, 3).Value = "Bíobío"
    $ws.Cells.Item($r, 4).Value = $newDate
    $ws.Cells.Item($r, 5).Value = 8
    $ws.Cells.Item($r, 6).Value = 100112027
    $ws.Cells.Item($r, 7).Value = "Melón"
    $ws.Cells.Item($r, 8).Value = $data[0]
    $ws.Cells.Item($r, 9).Value = $data[1]
    $ws.Cells.Item($r, 10).Value = $data[2]
    $ws.Cells.Item($r, 11).Value = $data[3]
    $ws.Cells.Item($r, 12).Value = $data[4]
    $ws.Cells.Item($r, 13).Value = $data[5]
    $ws.Cells.Item($r, 14).Value = "`$/unidad"
    $ws.Cells.Item($r, 15).Value = $data[6]
    $ws.Cells.Item($r, 16).Value = $data[7]
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
